$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A78").Value = "2025/12/06 00:00"
$ws.Range("B78").Value = "-"
$ws.Range("C78").Value = "-"
$ws.Range("D78").Value = "-"
$ws.Range("E78").Value = "-"
$ws.Range("F78").Value = "-"
$ws.Range("G78").Value = "-"
